# BUG: read_excel forward-filling MI names
#
# Adds two regression-test worksheets to the testmultiindex workbook:
#   - "mi_index_blank_after_name"      (copy of "mi_index_name", with the
#                                        cell right below the index name
#                                        blanked out instead of repeating it)
#   - "both_name_blank_after_mi_name"  (copy of "both_name", with the same
#                                        kind of cell blanked out)
#
# Each new sheet is inserted immediately after the sheet it was copied
# from, which is exactly where Excel's Worksheet.Copy(Before, After) puts
# it - matching the final tab order:
#   ... mi_index_name, mi_index_blank_after_name,
#       both_name, both_name_blank_after_mi_name,
#       both_name_skiprows, index_col_none

$wb = $excel.ActiveWorkbook

# --- mi_index_blank_after_name -------------------------------------------
$miIndexName = $wb.Worksheets.Item("mi_index_name")
$miIndexName.Copy($null, $miIndexName)
$miIndexBlank = $wb.Worksheets.Item("mi_index_name (2)")
$miIndexBlank.Name = "mi_index_blank_after_name"
# The row right under the repeated index name ("foo") used to repeat the
# name again in column B ("a"); it should be blank instead.
$miIndexBlank.Range("B2").ClearContents()

# --- both_name_blank_after_mi_name ---------------------------------------
$bothName = $wb.Worksheets.Item("both_name")
$bothName.Copy($null, $bothName)
$bothNameBlank = $wb.Worksheets.Item("both_name (2)")
$bothNameBlank.Name = "both_name_blank_after_mi_name"
$bothNameBlank.Range("B4").ClearContents()
